# "Added New Mac-Address and Document Types"
# Append one new data row (row 33) to the single worksheet, mirroring the
# existing rows' pattern (regcntr_id, usr_id, lang_code, is_active, cr_by,
# cr_dtimes), then move the selection the way the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item(1)

$newRow = 33

$ws.Cells.Item($newRow, 1).Value = 10002
$ws.Cells.Item($newRow, 2).Value = 110032
$ws.Cells.Item($newRow, 3).Value = "eng"
$ws.Cells.Item($newRow, 4).Value = $true
$ws.Cells.Item($newRow, 5).Value = "superadmin"
$ws.Cells.Item($newRow, 6).Value = "now()"

# Match the author's final on-screen selection / scroll position.
$ws.Range("C31").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
